$d = $word.ActiveDocument

function Replace-One($range, $old, $new) {
    $range.Find.ClearFormatting()
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 1)
}

# 1. Main body: "A QWREW," -> "A QWR,"
Replace-One $d.Content "QWREW" "QWR"

# 2. Header text replacements
$hdr = $d.Sections(1).Headers(1).Range

Replace-One $hdr "QWREW" "QWR"
Replace-One $hdr "REW" "QWER"
for ($i = 0; $i -lt 5; $i++) {
    Replace-One $hdr "Rew" "Qwer"
}
for ($i = 0; $i -lt 3; $i++) {
    Replace-One $hdr "rew" "qwer"
}
